$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: header "Save", styled like the other header cells (copy G1's format)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the Save values for rows 2-8
$saveValues = @(1, 0, 1, 0, 0, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}

$excel.CutCopyMode = 0
